$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.Shapes.Item(1).TextFrame.TextRange.Text = "Chapter 3"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Designing and Developing an Agent-based`rModel"
